$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -5968
# Row 98
$ws.Range("H98").Value = 4824.25
$ws.Range("I98").Value = 5255.909
$ws.Range("J98").Value = 76
$ws.Range("K98").Value = 5255.909
$ws.Range("L98").Value = 76
$ws.Range("M98").Value = -3757.909
$ws.Range("N98").Value = -3072
# Row 112
$ws.Range("H112").Value = 2802.5938
$ws.Range("J112").Value = 2919.4333
$ws.Range("L112").Value = 8758.2999
$ws.Range("N112").Value = -10974.2999
# Row 122
$ws.Range("H122").Value = 4824.25
$ws.Range("I122").Value = 5255.909
$ws.Range("J122").Value = 76
$ws.Range("K122").Value = 15767.727
$ws.Range("L122").Value = 228
$ws.Range("M122").Value = -13317.727
$ws.Range("N122").Value = -5128
# Row 132
$ws.Range("H132").Value = 15159908
$ws.Range("I132").Value = 16675520
$ws.Range("K132").Value = 50026560
$ws.Range("M132").Value = -50024030
# Row 135
$ws.Range("H135").Value = 602.3889
$ws.Range("I135").Value = 289.53333
$ws.Range("K135").Value = 2605.79997
$ws.Range("M135").Value = -70.79997000000003
# Row 138
$ws.Range("H138").Value = 1353.0613
$ws.Range("I138").Value = 869.2778
$ws.Range("J138").Value = 1633.9678
$ws.Range("K138").Value = 2607.8334
$ws.Range("L138").Value = 4901.903399999999
$ws.Range("M138").Value = 2532.1666
$ws.Range("N138").Value = -15181.9034
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1011.4828
$ws.Range("I110").Value = 588.6111
$ws.Range("K110").Value = 588.6111
$ws.Range("M110").Value = 1456.3889

$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 432.8421
$ws.Range("I80").Value = 307.85715
$ws.Range("K80").Value = 307.85715
$ws.Range("M80").Value = 690.14285
# Row 83
$ws.Range("H83").Value = 432.8421
$ws.Range("I83").Value = 307.85715
$ws.Range("K83").Value = 1539.28575
$ws.Range("M83").Value = 3452.71425
# Row 94
$ws.Range("H94").Value = 13889636
$ws.Range("I94").Value = 16667330
$ws.Range("J94").Value = 1166.3334
$ws.Range("K94").Value = 16667330
$ws.Range("L94").Value = 1166.3334
$ws.Range("M94").Value = -16666879
$ws.Range("N94").Value = -2068.3334
# Row 99
$ws.Range("H99").Value = 31251406
$ws.Range("I99").Value = 35715656
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 35715656
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -35714158
$ws.Range("N99").Value = -4646
# Row 134
$ws.Range("H134").Value = 7271.591
$ws.Range("I134").Value = 1131.7333
$ws.Range("K134").Value = 3395.199900000001
$ws.Range("M134").Value = -860.1999000000005

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 34483540
$ws.Range("I16").Value = 38462256
$ws.Range("J16").Value = 1320
$ws.Range("K16").Value = 38462256
$ws.Range("L16").Value = 1320
$ws.Range("M16").Value = -38461969
$ws.Range("N16").Value = -1894
# Row 58
$ws.Range("H58").Value = 918.87933
$ws.Range("I58").Value = 723.1627999999999
$ws.Range("K58").Value = 723.1627999999999
$ws.Range("M58").Value = -520.1627999999999
# Row 62
$ws.Range("H62").Value = 66670000
$ws.Range("I62").Value = 4999.5
$ws.Range("J62").Value = 200000000
$ws.Range("K62").Value = 4999.5
$ws.Range("L62").Value = 200000000
$ws.Range("M62").Value = -4375.5
$ws.Range("N62").Value = -200001248
# Row 65
$ws.Range("H65").Value = 66670000
$ws.Range("I65").Value = 4999.5
$ws.Range("J65").Value = 200000000
$ws.Range("K65").Value = 24997.5
$ws.Range("L65").Value = 1000000000
$ws.Range("M65").Value = -21877.5
$ws.Range("N65").Value = -1000006240
# Row 105
$ws.Range("H105").Value = 739.8
$ws.Range("I105").Value = 674.75
$ws.Range("K105").Value = 674.75
$ws.Range("M105").Value = 1072.25
# Row 113
$ws.Range("H113").Value = 34483540
$ws.Range("I113").Value = 38462256
$ws.Range("J113").Value = 1320
$ws.Range("K113").Value = 38462256
$ws.Range("L113").Value = 1320
$ws.Range("M113").Value = -38460086
$ws.Range("N113").Value = -5660
# Row 122
$ws.Range("H122").Value = 956.1905
$ws.Range("I122").Value = 795.4
$ws.Range("K122").Value = 2386.2
$ws.Range("M122").Value = 63.80000000000018
# Row 136
$ws.Range("H136").Value = 918.87933
$ws.Range("I136").Value = 723.1627999999999
$ws.Range("K136").Value = 2169.4884
$ws.Range("M136").Value = 380.5116000000003

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 763.4091
$ws.Range("I107").Value = 828.0769
$ws.Range("J107").Value = 670
$ws.Range("K107").Value = 828.0769
$ws.Range("L107").Value = 670
$ws.Range("M107").Value = 1091.9231
$ws.Range("N107").Value = -4510
# Row 113
$ws.Range("H113").Value = 889.875
$ws.Range("I113").Value = 874.1429000000001
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 874.1429000000001
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 1295.8571
$ws.Range("N113").Value = -5340
# Row 122
$ws.Range("H122").Value = 1958.9445
$ws.Range("I122").Value = 1607.9286
$ws.Range("J122").Value = 3187.5
$ws.Range("K122").Value = 4823.7858
$ws.Range("L122").Value = 9562.5
$ws.Range("M122").Value = -2373.7858
$ws.Range("N122").Value = -14462.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1412.5
$ws.Range("I7").Value = 1498.3334
$ws.Range("K7").Value = 1498.3334
$ws.Range("M7").Value = -1386.3334
# Row 46
$ws.Range("H46").Value = 2999.5
$ws.Range("J46").Value = 2999.5
$ws.Range("L46").Value = 2999.5
$ws.Range("N46").Value = -3375.5
# Row 126
$ws.Range("H126").Value = 1412.5
$ws.Range("I126").Value = 1498.3334
$ws.Range("K126").Value = 4495.0002
$ws.Range("M126").Value = -2025.0002

$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 16589.5
$ws.Range("J52").Value = 16589.5
$ws.Range("L52").Value = 16589.5
$ws.Range("N52").Value = -17041.5
# Row 100
$ws.Range("H100").Value = 351.7143
$ws.Range("I100").Value = 310.33334
$ws.Range("K100").Value = 620.66668
$ws.Range("M100").Value = -79.66668000000004
# Row 107
$ws.Range("H107").Value = 340.0909
$ws.Range("I107").Value = 306.66666
$ws.Range("J107").Value = 380.2
$ws.Range("K107").Value = 919.9999799999999
$ws.Range("L107").Value = 1140.6
$ws.Range("M107").Value = 1000.00002
$ws.Range("N107").Value = -4980.6
# Row 122
$ws.Range("H122").Value = 33159678
$ws.Range("I122").Value = 37060590
$ws.Range("K122").Value = 111181770
$ws.Range("M122").Value = -111179320
